$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (family_ablity) values for rows 2-11 ---
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 0

# --- gender column (F) updates ---
$ws.Range("F2").Value = "M"
$ws.Range("F7").Value = "F"

# --- New row 12 : 독수리 (eagle) ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "독수리"
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "M"

# --- Column D width ---
$ws.Range("D1").ColumnWidth = 13.714285714285714

# --- Page setup (printer settings) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("J8").Select()
